# edit.ps1
# Applies the "Yais SZE_timetable" updates:
#  - Adjust a few column widths
#  - Update several activity labels (name formatting, master class entries)
#  - Add a "Free Time" entry in J24 (splitting what used to be a single J20:J27 block)
#  - Remove the "Free Time" row that existed at 17:15 (row 32) for days 1-5
#  - Re-shape the row 28-31 merged blocks to stop at row 30 instead of row 31,
#    and remove the now-unused row 32-39 merged blocks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column width changes -------------------------------------------------
# COM ColumnWidth is reported ~0.8333 (5/6) wider than the character-width value
# stored in the OOXML <col width="..."/> attribute, so we subtract that offset
# to land exactly on the target stored width.
$widthOffset = 5 / 6
$ws.Range("B1").EntireColumn.ColumnWidth = (35 - $widthOffset)
$ws.Range("D1").EntireColumn.ColumnWidth = (32 - $widthOffset)
$ws.Range("F1").EntireColumn.ColumnWidth = (35 - $widthOffset)
$ws.Range("J1").EntireColumn.ColumnWidth = (34 - $widthOffset)

# --- 2. Unmerge all the cell blocks we need to touch (merged ranges only
#        allow writes through their top-left anchor cell) ---------------------
$ws.Range("J20:J27").UnMerge()
$ws.Range("J28:J31").UnMerge()
$ws.Range("B28:B31").UnMerge()
$ws.Range("D28:D31").UnMerge()
$ws.Range("F28:F31").UnMerge()
$ws.Range("H28:H31").UnMerge()
$ws.Range("B32:B39").UnMerge()
$ws.Range("D32:D39").UnMerge()
$ws.Range("F32:F39").UnMerge()
$ws.Range("H32:H39").UnMerge()
$ws.Range("J32:J39").UnMerge()

# --- 3. Text updates -----------------------------------------------------
$ws.Range("D7").Value = "Private lesson with Ivy CHUANG"
$ws.Range("J7").Value = "Master class with Ivy & Stephane"
$ws.Range("B11").Value = "Private lesson with Stephane RETY"
$ws.Range("F20").Value = "Private lesson with Stephane RETY"
$ws.Range("J20").Value = "Master class with Ivy & Stephane"

# New activity cell for Day 5 at 15:15 (row 24)
$ws.Range("J24").Value = "Free Time"
$ws.Range("J24").VerticalAlignment = -4108  # xlCenter, matches style used by the other activity cells

# --- 4. Remove the 17:15 (row 32) "Free Time" entries for Days 1-5 -----------
foreach ($col in @("B", "D", "F", "H", "J")) {
    $cell = $ws.Range($col + "32")
    $cell.Value = ""
    $cell.ClearFormats()
}

# --- 5. Recreate merges with the new shapes ----------------------------------
$ws.Range("B28:B30").Merge()
$ws.Range("D28:D30").Merge()
$ws.Range("F28:F30").Merge()
$ws.Range("H28:H30").Merge()
$ws.Range("J20:J23").Merge()
$ws.Range("J24:J27").Merge()
$ws.Range("J28:J30").Merge()

# --- 6. Clean up stray auto-populated cells ----------------------------------
# Merging a range through this COM layer stamps every cell in the range with
# the anchor cell's style, even for rows that previously had no data in that
# column. Blank those back out (value + format) so the saved sheet only has
# real content where the source workbook has it.
$strayCells = @(
    "J21", "J22", "J23",
    "J25", "J26", "J27",
    "B29", "D29", "F29", "H29", "J29",
    "B30", "D30", "F30", "H30", "J30"
)
foreach ($addr in $strayCells) {
    $cell = $ws.Range($addr)
    $cell.Value = ""
    $cell.ClearFormats()
}
